# Apply "updated survay data format" changes:
#  1. Rename sheet from "survay_data_test" to "survay_data"
#  2. Replace "NA" with "NULL" in a specific set of cells (adds a new shared string)
#  3. Move active selection to A4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet.
$ws.Name = "survay_data"

# 2. Update specific cells from "NA" to "NULL".
$cellsToUpdate = @("E2", "E4", "E5", "D28", "E28", "D32", "E32", "D34", "E34", "E39", "D40", "E40", "F40", "E46")
foreach ($cellRef in $cellsToUpdate) {
    $ws.Range($cellRef).Value = "NULL"
}

# 3. Update the active selection to A4.
$ws.Range("A4").Select()
